$d = $word.ActiveDocument

# --- Change 1: expand the "good model" sentence into the "decent model" critique. ---
$old1 = "Based on the accuracy and that the terms are all significant, this is a good model."
$new1 = "Based on the accuracy and that the terms are all significant, this is a decent model. From the contingency table, there are relatively few loans predicted as bad compared to good loans. The model may be good at predicting good loans, and bad at predicting bad ones."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# --- Change 2: the accuracy plot's paragraph accidentally trailed off into a stray
#     "## 7. Optimizing the Threshold for Profit" heading fragment; strip that text so the
#     paragraph holds only the picture. ---
$d.Content.Find.Execute(" ## 7. Optimizing the Threshold for Profit", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# Locate the (now picture-only) accuracy-plot paragraph and the old
# "## [1] 0.7586207" SourceCode paragraph that used to follow the stray heading text.
$accFigIdx = -1
$srcIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*0.7586207*") {
        $srcIdx = $i
    }
}
$accFigIdx = $srcIdx - 1

# --- Change 3: add the new "maximum accuracy" commentary paragraph right after the
#     accuracy plot. ---
$accPara = $d.Paragraphs.Item($accFigIdx)
$accPara.Range.InsertParagraphAfter()
$newBodyPara = $d.Paragraphs.Item($accFigIdx + 1)
$newBodyPara.Style = "BodyText"
$newBodyPara.Range.Text = "The maxmimum accuracy occurs near the threshold of .55. As you can see from the plot above of accuracy vs threshold, the accuracy is mostly flat from 0 to 0.6, but from inspection of the data peaks at .55. It then declines between 0.6 and 1.0."

# --- Change 4: add the real "7. Optimizing the Threshold for Profit" Heading2 (with its
#     bookmark) right after that new paragraph, then drop the stale SourceCode paragraph
#     that used to hold "## [1] 0.7586207". ---
$newBodyPara.Range.InsertParagraphAfter()
$newHeadingPara = $d.Paragraphs.Item($accFigIdx + 2)
$newHeadingPara.Style = "Heading2"
$newHeadingPara.Range.Text = "7. Optimizing the Threshold for Profit"
$collapsed = $d.Range($newHeadingPara.Range.Start, $newHeadingPara.Range.Start)
$d.Bookmarks.Add("optimizing-the-threshold-for-profit", $collapsed) | Out-Null

$srcIdx2 = $srcIdx + 2
$srcPara = $d.Paragraphs.Item($srcIdx2)
$srcPara.Range.Delete()

# --- Change 5: the profit plot's picture paragraph turns into a text paragraph describing
#     the maximum-profit threshold, followed by two new BodyText paragraphs. ---
$profitIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.InlineShapes.Count -gt 0) {
        $profitIdx = $i
    }
}
$profitPara = $d.Paragraphs.Item($profitIdx)
$profitPara.Range.Text = "The maximum profit of `$2,609,607 occurs at a threshold of .76. The best profit threshold corresponds to an accuracy of 65%. The profit where the model is the most accurate at a threshold of 55% is `$1,247,186. The maximum accuracy and profit thresholds do not coincide."

$profitPara.Range.InsertParagraphAfter()
$bodyPara1 = $d.Paragraphs.Item($profitIdx + 1)
$bodyPara1.Style = "BodyText"
$bodyPara1.Range.Text = "If no loans were denied the total profit would be `$758,390. The increase from using the model is 244%."

$bodyPara1.Range.InsertParagraphAfter()
$bodyPara2 = $d.Paragraphs.Item($profitIdx + 2)
$bodyPara2.Style = "BodyText"
$bodyPara2.Range.Text = "If the model had predicted all the good loans perfectly, the total profit would be `$12,020,435, which represents an increase of 1485%. Perfect prediction is at least 6 times better, so may be the model could be improved."

# --- Change 6: add the closing takeaway paragraph right after "8. Results Summary". ---
$resultsIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*8. Results Summary*") {
        $resultsIdx = $i
    }
}
$resultsPara = $d.Paragraphs.Item($resultsIdx)
$resultsPara.Range.InsertParagraphAfter()
$finalPara = $d.Paragraphs.Item($resultsIdx + 1)
$finalPara.Style = "FirstParagraph"
$finalPara.Range.Text = "The classification threshold that produces the maximum profit is .76, which gives an accuracy of 65%."
